$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "Category 1" -> "testupdateonform" (shared string used by B2)
$ws.Range("B2").Value = "testupdateonform"

# Update row 4 (Category ID 3 -> 5, keep Name/Description "newcatetest")
$ws.Range("A4").Value = 5

# Update row 5 (was Category ID 4 / newcatetest / newcatetest -> now 6 / test / test)
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "test"
$ws.Range("C5").Value = "test"

# Remove old rows 6, 7 and 8 entirely (data previously held by ids 5,6,7)
$ws.Range("A6:C8").Delete()
